$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text/value updates derived from the authoritative OOXML diff.
# Cells whose new value is numeric-looking (e.g. "97.80") are forced to
# Text format first so Excel keeps the exact original string (trailing
# zeros, dot-grouped thousands separators, etc.) instead of silently
# re-parsing it into a Number and losing formatting/precision.

$ws.Range('D2').Value = '42.824.81'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.540.27'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.04'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.80'
$ws.Range('E6').Value = '  +5.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.79'
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0825'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.114'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.928.22'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '2.503.93'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.06'
$ws.Range('E16').Value = '  +5.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.871'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '42.829.12'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.29'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('D20').Value = '0.0₃0990'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.56'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.63'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.01'
$ws.Range('E23').Value = '  -0.55%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.08'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.94'
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.18'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  +6.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.05'
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.19'
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.25'
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.60'
$ws.Range('E33').Value = '  +13.84%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.31'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0796'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.115'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.04'
$ws.Range('E39').Value = '  +5.33%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  +30.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.87'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').Value = '2.090.78'
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0305'
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.60'
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.94'
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('D49').Value = '2.786.27'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.68'
$ws.Range('E50').Value = '  +6.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.192'
$ws.Range('E51').Value = '  +1.70%  '
